$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "size" (column B) and "size_design" (column C) values from
# 300/200 (ms) to 0.3/0.2 (s) for rows 2-7 so the drag-and-drop demo
# doesn't end on the first click.
$ws.Range("B2:B7").Value = 0.3
$ws.Range("C2:C7").Value = 0.2

# Update the active cell selection to match the saved state in the diff.
$ws.Range("C7").Select()
